$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 13. This shifts the existing
# row 13 -> row 14 and row 14 -> row 15, preserving their data intact.
$ws.Rows.Item(13).Insert()

# Fill the new row 13 with a fresh weekly price report for
# Vega Monumental Concepción - Arveja Verde.
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Vega Monumental Concepción"
$ws.Range("C13").Value = "Bíobío"
$ws.Range("D13").Value = 44518
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 100112022
$ws.Range("G13").Value = "Arveja Verde"
$ws.Range("H13").Value = "Perfection"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 350
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14571
$ws.Range("N13").Value = "$/saco 25 kilos"
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 583
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
